$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the case number in A2 (the target-path value)
$ws.Range("A2").Value = "250100032HZH"

# Move the selection to B3, matching the recorded cursor position
$ws.Range("B3").Select()
